# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.17 = 20589.15 pesos`n✅ 20589.15 pesos = 5.14 = 953.63 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas: update the updated exchange rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 193.5
$wsTasas.Range("O10").Value = 3984

$wsTasas.Range("N12").Value = 4004.99
$wsTasas.Range("O12").Value = 185.5
